# Auto-generated edit script: updates market-data value cells
# per the authoritative diff (scheduled runner data refresh).
$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(6, 8).Value = 2652668.8
$ws.Cells.Item(6, 9).Value = 2652668.8
$ws.Cells.Item(6, 11).Value = 7958006.399999999
$ws.Cells.Item(6, 13).Value = -7957894.399999999
$ws.Cells.Item(8, 8).Value = 312.64444
$ws.Cells.Item(8, 9).Value = 352.75
$ws.Cells.Item(8, 11).Value = 1058.25
$ws.Cells.Item(8, 13).Value = -919.25
$ws.Cells.Item(86, 8).Value = 115387450
$ws.Cells.Item(86, 10).Value = 18524518
$ws.Cells.Item(86, 12).Value = 18524518
$ws.Cells.Item(86, 14).Value = -18526764
$ws.Cells.Item(89, 8).Value = 115387450
$ws.Cells.Item(89, 10).Value = 18524518
$ws.Cells.Item(89, 12).Value = 92622590
$ws.Cells.Item(89, 14).Value = -92633822
$ws.Cells.Item(132, 8).Value = 1444.66
$ws.Cells.Item(132, 9).Value = 1192.1086
$ws.Cells.Item(132, 10).Value = 4349
$ws.Cells.Item(132, 11).Value = 3576.3258
$ws.Cells.Item(132, 12).Value = 13047
$ws.Cells.Item(132, 13).Value = -1046.3258
$ws.Cells.Item(132, 14).Value = -18107
$ws.Cells.Item(136, 8).Value = 53593.332
$ws.Cells.Item(136, 10).Value = 60780
$ws.Cells.Item(136, 12).Value = 60780
$ws.Cells.Item(136, 14).Value = -70980
$ws.Cells.Item(137, 8).Value = 2791.6875
$ws.Cells.Item(137, 9).Value = 3063.3
$ws.Cells.Item(137, 11).Value = 9189.900000000001
$ws.Cells.Item(137, 13).Value = -6639.900000000001
# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 13).ClearContents()
$ws.Cells.Item(61, 8).Value = 4899.4375
$ws.Cells.Item(61, 9).Value = 2400.7878
$ws.Cells.Item(61, 11).Value = 2400.7878
$ws.Cells.Item(61, 13).Value = -2188.7878
$ws.Cells.Item(74, 8).Value = 18873.684
$ws.Cells.Item(74, 9).Value = 24044.822
$ws.Cells.Item(74, 11).Value = 24044.822
$ws.Cells.Item(74, 13).Value = -23170.822
$ws.Cells.Item(77, 8).Value = 18873.684
$ws.Cells.Item(77, 9).Value = 24044.822
$ws.Cells.Item(77, 11).Value = 120224.11
$ws.Cells.Item(77, 13).Value = -115856.11
$ws.Cells.Item(110, 8).Value = 27778964
$ws.Cells.Item(110, 9).Value = 1292.3636
$ws.Cells.Item(110, 11).Value = 1292.3636
$ws.Cells.Item(110, 13).Value = 752.6364000000001
$ws.Cells.Item(132, 8).Value = 8099.8066
$ws.Cells.Item(132, 9).Value = 6573.9473
$ws.Cells.Item(132, 11).Value = 19721.8419
$ws.Cells.Item(132, 13).Value = -17191.8419
$ws.Cells.Item(136, 8).Value = 4899.4375
$ws.Cells.Item(136, 9).Value = 2400.7878
$ws.Cells.Item(136, 11).Value = 7202.3634
$ws.Cells.Item(136, 13).Value = -4652.3634
# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(86, 8).Value = 83403730
$ws.Cells.Item(86, 9).Value = 31377550
$ws.Cells.Item(86, 11).Value = 31377550
$ws.Cells.Item(86, 13).Value = -31376427
$ws.Cells.Item(89, 8).Value = 83403730
$ws.Cells.Item(89, 9).Value = 31377550
$ws.Cells.Item(89, 11).Value = 156887750
$ws.Cells.Item(89, 13).Value = -156882134
$ws.Cells.Item(94, 8).Value = 1710.129
$ws.Cells.Item(94, 9).Value = 530.6
$ws.Cells.Item(94, 11).Value = 530.6
$ws.Cells.Item(94, 13).Value = -79.60000000000002
$ws.Cells.Item(105, 8).Value = 3340.8
$ws.Cells.Item(105, 9).Value = 2929.7646
$ws.Cells.Item(105, 10).Value = 4214.25
$ws.Cells.Item(105, 11).Value = 2929.7646
$ws.Cells.Item(105, 12).Value = 4214.25
$ws.Cells.Item(105, 13).Value = -1182.7646
$ws.Cells.Item(105, 14).Value = -7708.25
$ws.Cells.Item(134, 8).Value = 7320.7354
$ws.Cells.Item(134, 9).Value = 4150
$ws.Cells.Item(134, 11).Value = 12450
$ws.Cells.Item(134, 13).Value = -9915
# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(22, 8).Value = 469.1
$ws.Cells.Item(22, 9).Value = 449.25
$ws.Cells.Item(22, 10).Value = 548.5
$ws.Cells.Item(22, 11).Value = 449.25
$ws.Cells.Item(22, 12).Value = 548.5
$ws.Cells.Item(22, 13).Value = -99.25
$ws.Cells.Item(22, 14).Value = -1248.5
$ws.Cells.Item(35, 8).Value = 15038.25
$ws.Cells.Item(35, 9).Value = 15050
$ws.Cells.Item(35, 10).Value = 15026.5
$ws.Cells.Item(35, 11).Value = 15050
$ws.Cells.Item(35, 12).Value = 15026.5
$ws.Cells.Item(35, 13).Value = -14756
$ws.Cells.Item(35, 14).Value = -15614.5
$ws.Cells.Item(62, 8).Value = 17863858
$ws.Cells.Item(62, 10).Value = 8000.4
$ws.Cells.Item(62, 12).Value = 8000.4
$ws.Cells.Item(62, 14).Value = -9248.4
$ws.Cells.Item(65, 8).Value = 17863858
$ws.Cells.Item(65, 10).Value = 8000.4
$ws.Cells.Item(65, 12).Value = 40002
$ws.Cells.Item(65, 14).Value = -46242
$ws.Cells.Item(132, 8).Value = 7157.913
$ws.Cells.Item(132, 9).Value = 2617.8572
$ws.Cells.Item(132, 11).Value = 7853.571599999999
$ws.Cells.Item(132, 13).Value = -5323.571599999999
$ws.Cells.Item(134, 8).Value = 7436.9
$ws.Cells.Item(134, 9).Value = 6916.4194
$ws.Cells.Item(134, 11).Value = 20749.2582
$ws.Cells.Item(134, 13).Value = -18214.2582
# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(5, 8).Value = 3640450
$ws.Cells.Item(5, 9).Value = 6668258.5
$ws.Cells.Item(5, 10).Value = 7079.8
$ws.Cells.Item(5, 11).Value = 20004775.5
$ws.Cells.Item(5, 12).Value = 21239.4
$ws.Cells.Item(5, 13).Value = -20004663.5
$ws.Cells.Item(5, 14).Value = -21463.4
$ws.Cells.Item(7, 8).Value = 1740
$ws.Cells.Item(7, 9).Value = 266.66666
$ws.Cells.Item(7, 10).Value = 2371.4285
$ws.Cells.Item(7, 11).Value = 799.9999799999999
$ws.Cells.Item(7, 12).Value = 7114.2855
$ws.Cells.Item(7, 13).Value = -687.9999799999999
$ws.Cells.Item(7, 14).Value = -7338.2855
$ws.Cells.Item(88, 8).Value = 2000
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 14).ClearContents()
$ws.Cells.Item(91, 8).Value = 2000
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 5768.222
$ws.Cells.Item(113, 10).Value = 6158.0625
$ws.Cells.Item(113, 12).Value = 18474.1875
$ws.Cells.Item(113, 14).Value = -22814.1875
$ws.Cells.Item(131, 8).Value = 1995.8438
$ws.Cells.Item(131, 10).Value = 2191.7917
$ws.Cells.Item(131, 12).Value = 6575.375100000001
$ws.Cells.Item(131, 14).Value = -16655.3751
$ws.Cells.Item(135, 8).Value = 3640450
$ws.Cells.Item(135, 9).Value = 6668258.5
$ws.Cells.Item(135, 10).Value = 7079.8
$ws.Cells.Item(135, 11).Value = 60014326.5
$ws.Cells.Item(135, 12).Value = 63718.2
$ws.Cells.Item(135, 13).Value = -60011791.5
$ws.Cells.Item(135, 14).Value = -68788.20000000001
# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(80, 8).Value = 2871.9333
$ws.Cells.Item(80, 9).Value = 2632.5557
$ws.Cells.Item(80, 10).Value = 3231
$ws.Cells.Item(80, 11).Value = 2632.5557
$ws.Cells.Item(80, 12).Value = 3231
$ws.Cells.Item(80, 13).Value = -1634.5557
$ws.Cells.Item(80, 14).Value = -5227
$ws.Cells.Item(83, 8).Value = 2871.9333
$ws.Cells.Item(83, 9).Value = 2632.5557
$ws.Cells.Item(83, 10).Value = 3231
$ws.Cells.Item(83, 11).Value = 13162.7785
$ws.Cells.Item(83, 12).Value = 16155
$ws.Cells.Item(83, 13).Value = -8170.7785
$ws.Cells.Item(83, 14).Value = -26139
$ws.Cells.Item(113, 8).Value = 5141.36
$ws.Cells.Item(113, 9).Value = 2323.818
$ws.Cells.Item(113, 10).Value = 7355.143
$ws.Cells.Item(113, 11).Value = 2323.818
$ws.Cells.Item(113, 12).Value = 7355.143
$ws.Cells.Item(113, 13).Value = -153.8180000000002
$ws.Cells.Item(113, 14).Value = -11695.143
$ws.Cells.Item(132, 8).Value = 15577
$ws.Cells.Item(132, 9).Value = 9000
$ws.Cells.Item(132, 10).Value = 16673.166
$ws.Cells.Item(132, 11).Value = 27000
$ws.Cells.Item(132, 12).Value = 50019.49800000001
$ws.Cells.Item(132, 13).Value = -24470
$ws.Cells.Item(132, 14).Value = -55079.49800000001
# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(9, 8).Value = 640
$ws.Cells.Item(9, 9).Value = 640
$ws.Cells.Item(9, 11).Value = 640
$ws.Cells.Item(9, 13).Value = -416
$ws.Cells.Item(16, 8).Value = 1458.138
$ws.Cells.Item(16, 9).Value = 1403.1072
$ws.Cells.Item(16, 11).Value = 1403.1072
$ws.Cells.Item(16, 13).Value = -1233.1072
$ws.Cells.Item(46, 8).Value = 7409749.5
$ws.Cells.Item(46, 9).Value = 899
$ws.Cells.Item(46, 11).Value = 899
$ws.Cells.Item(46, 13).Value = -711
$ws.Cells.Item(82, 8).Value = 1085553.9
$ws.Cells.Item(82, 9).Value = 1567000.6
$ws.Cells.Item(82, 11).Value = 1567000.6
$ws.Cells.Item(82, 13).Value = -1566639.6
$ws.Cells.Item(85, 8).Value = 1085553.9
$ws.Cells.Item(85, 9).Value = 1567000.6
$ws.Cells.Item(85, 11).Value = 1567000.6
$ws.Cells.Item(85, 13).Value = -1565752.6
$ws.Cells.Item(93, 8).Value = 4858.45
$ws.Cells.Item(93, 9).Value = 4739.067
$ws.Cells.Item(93, 11).Value = 4739.067
$ws.Cells.Item(93, 13).Value = -3491.067
$ws.Cells.Item(136, 8).Value = 10719.652
$ws.Cells.Item(136, 9).Value = 2365.1428
$ws.Cells.Item(136, 11).Value = 7095.428400000001
$ws.Cells.Item(136, 13).Value = -4545.428400000001
# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(132, 8).Value = 9268052
$ws.Cells.Item(132, 9).Value = 12503923
$ws.Cells.Item(132, 11).Value = 37511769
$ws.Cells.Item(132, 13).Value = -37509239
